# Apply updated cosinor statistics (CircaDB / CircadiPy re-run) to rows 2-15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = "[0.030424206341297833, 14.775346773624708]"
$ws.Range("N2").Value = 0.04910122339608458
$ws.Range("O2").Value = 0.04910122339608458
$ws.Range("Q2").Value = "[-2.7925268031909276, 0.12578949563923203]"
$ws.Range("R2").Value = 0.07229865087324039
$ws.Range("S2").Value = 0.07229865087324039
$ws.Range("U2").Value = "[5.383008004088298, 13.649291308406374]"
$ws.Range("V2").Value = [double]"3.04530778068024e-05"
$ws.Range("W2").Value = [double]"3.04530778068024e-05"
$ws.Range("Y2").Value = -0.4554554554554633
$ws.Range("Z2").Value = 10.11111111111116

# Row 3
$ws.Range("M3").Value = "[-0.42455168178672764, 15.37476150240607]"
$ws.Range("N3").Value = 0.06306835207361994
$ws.Range("O3").Value = 0.06306835207361994
$ws.Range("Q3").Value = "[-2.943474197958005, 0.8553685703467702]"
$ws.Range("R3").Value = 0.2741380081293789
$ws.Range("S3").Value = 0.2741380081293789
$ws.Range("U3").Value = "[6.303862568273926, 14.867295647033204]"
$ws.Range("V3").Value = [double]"9.851388260262439e-06"
$ws.Range("W3").Value = [double]"9.851388260262439e-06"
$ws.Range("Y3").Value = -3.097097097097111
$ws.Range("Z3").Value = 10.65765765765772

# Row 4
$ws.Range("M4").Value = "[-0.20841625806219177, 13.961323182813587]"
$ws.Range("N4").Value = 0.05682924728112781
$ws.Range("O4").Value = 0.05682924728112781
$ws.Range("Q4").Value = "[-2.7170531058073886, 0.8302106712189241]"
$ws.Range("R4").Value = 0.2897327401351861
$ws.Range("S4").Value = 0.2897327401351861
$ws.Range("U4").Value = "[5.01424413443368, 12.591822464212571]"
$ws.Range("V4").Value = [double]"2.651021426181011e-05"
$ws.Range("W4").Value = [double]"2.651021426181011e-05"
$ws.Range("Y4").Value = -3.00600600600602
$ws.Range("Z4").Value = 9.83783783783789

# Row 5
$ws.Range("M5").Value = "[-1.1367741067723944, 12.661401120364571]"
$ws.Range("N5").Value = 0.09945066554806026
$ws.Range("O5").Value = 0.09945066554806026
$ws.Range("Q5").Value = "[-2.295658295415964, 1.540921321580579]"
$ws.Range("R5").Value = 0.6938174732721123
$ws.Range("S5").Value = 0.6938174732721123
$ws.Range("U5").Value = "[4.329777940187238, 12.010820210711056]"
$ws.Range("V5").Value = [double]"9.507838126787682e-05"
$ws.Range("W5").Value = [double]"9.507838126787682e-05"
$ws.Range("Y5").Value = -5.579329329329358
$ws.Range("Z5").Value = 8.3120620620621

# Row 6
$ws.Range("M6").Value = "[-1.889325657931165, 12.395663130582937]"
$ws.Range("N6").Value = 0.1454861762371933
$ws.Range("O6").Value = 0.1454861762371933
$ws.Range("Q6").Value = "[-2.9749215718678124, 2.9623426223038885]"
$ws.Range("R6").Value = 0.9966141639383523
$ws.Range("S6").Value = 0.9966141639383523
$ws.Range("U6").Value = "[4.791408376435093, 13.374788043200972]"
$ws.Range("V6").Value = 0.0001020033447041158
$ws.Range("W6").Value = 0.0001020033447041158
$ws.Range("Y6").Value = -10.72597597597603
$ws.Range("Z6").Value = 10.77152152152158

# Row 7
$ws.Range("M7").Value = "[-1.2881337877825434, 12.695694885346265]"
$ws.Range("N7").Value = 0.1073447283787701
$ws.Range("O7").Value = 0.1073447283787701
$ws.Range("Q7").Value = "[-1.547210796362541, 2.1887372241226197]"
$ws.Range("R7").Value = 0.731062589904528
$ws.Range("S7").Value = 0.731062589904528
$ws.Range("U7").Value = "[4.288891337900383, 12.00678414281378]"
$ws.Range("V7").Value = 0.0001053411395197656
$ws.Range("W7").Value = 0.0001053411395197656
$ws.Range("Y7").Value = 14.82507507507515
$ws.Range("Z7").Value = 28.35210210210225

# Row 8
$ws.Range("M8").Value = "[-1.4989631982729534, 14.658623843870185]"
$ws.Range("N8").Value = 0.1078945018456265
$ws.Range("O8").Value = 0.1078945018456265
$ws.Range("Q8").Value = "[-2.5535267614763884, 3.723369070921236]"
$ws.Range("R8").Value = 0.709145957153668
$ws.Range("S8").Value = 0.709145957153668
$ws.Range("U8").Value = "[5.580333430068341, 13.936355272871177]"
$ws.Range("V8").Value = [double]"2.445830853226028e-05"
$ws.Range("W8").Value = [double]"2.445830853226028e-05"
$ws.Range("Y8").Value = 9.268518518518565
$ws.Range("Z8").Value = 31.99574574574591

# Row 9
$ws.Range("M9").Value = "[0.3396345114172714, 14.961044228441073]"
$ws.Range("N9").Value = 0.04066259261658711
$ws.Range("O9").Value = 0.04066259261658711
$ws.Range("Q9").Value = "[-0.5912106295043857, 2.0126319302276956]"
$ws.Range("R9").Value = 0.277403426954332
$ws.Range("S9").Value = 0.277403426954332
$ws.Range("U9").Value = "[5.532314164306788, 13.297688285882797]"
$ws.Range("V9").Value = [double]"1.352525110531744e-05"
$ws.Range("W9").Value = [double]"1.352525110531744e-05"
$ws.Range("Y9").Value = 15.46271271271279
$ws.Range("Z9").Value = 24.89064064064077

# Row 10
$ws.Range("M10").Value = "[-1.145859456939224, 14.689845590475255]"
$ws.Range("N10").Value = 0.09182323194601105
$ws.Range("O10").Value = 0.09182323194601105
$ws.Range("Q10").Value = "[-0.5346053564667317, 2.937184723176043]"
$ws.Range("R10").Value = 0.1702180043819306
$ws.Range("S10").Value = 0.1702180043819306
$ws.Range("U10").Value = "[5.1743240941912845, 13.563491664723294]"
$ws.Range("V10").Value = [double]"4.780945606563947e-05"
$ws.Range("W10").Value = [double]"4.780945606563947e-05"
$ws.Range("Y10").Value = 12.92988988989008
$ws.Range("Z10").Value = 26.34586586586626

# Row 11
$ws.Range("M11").Value = "[-0.8558811548603718, 14.746073556519406]"
$ws.Range("N11").Value = 0.07967520619979696
$ws.Range("O11").Value = 0.07967520619979696
$ws.Range("Q11").Value = "[-0.7107106503616549, 2.96234262230389]"
$ws.Range("R11").Value = 0.2233606940443515
$ws.Range("S11").Value = 0.2233606940443515
$ws.Range("U11").Value = "[5.284196591920256, 13.758054394032108]"
$ws.Range("V11").Value = [double]"4.37458798778767e-05"
$ws.Range("W11").Value = [double]"4.37458798778767e-05"
$ws.Range("Y11").Value = 12.83267267267286
$ws.Range("Z11").Value = 27.02638638638678

# Row 12
$ws.Range("M12").Value = "[-0.39936005695823873, 14.543863342839414]"
$ws.Range("N12").Value = 0.0629899248479957
$ws.Range("O12").Value = 0.0629899248479957
$ws.Range("Q12").Value = "[-2.0503687789194647, 3.5724216761541605]"
$ws.Range("R12").Value = 0.5883049860406935
$ws.Range("S12").Value = 0.5883049860406935
$ws.Range("U12").Value = "[4.765801078815519, 12.53119161826656]"
$ws.Range("V12").Value = [double]"4.975804018858554e-05"
$ws.Range("W12").Value = [double]"4.975804018858554e-05"
$ws.Range("Y12").Value = 10.4751551551553
$ws.Range("Z12").Value = 32.20320320320367

# Row 13
$ws.Range("M13").Value = "[-0.19164733717850524, 14.737681168131097]"
$ws.Range("N13").Value = 0.05591792325550538
$ws.Range("O13").Value = 0.05591792325550538
$ws.Range("Q13").Value = "[-1.3333686537758478, 2.7044741562434655]"
$ws.Range("R13").Value = 0.4975346913099163
$ws.Range("S13").Value = 0.4975346913099163
$ws.Range("U13").Value = "[4.906548306151884, 12.60968373534595]"
$ws.Range("V13").Value = [double]"3.672196195747546e-05"
$ws.Range("W13").Value = [double]"3.672196195747546e-05"
$ws.Range("Y13").Value = 13.82914914914935
$ws.Range("Z13").Value = 29.43251251251294

# Row 14
$ws.Range("M14").Value = "[0.2255067007601692, 14.752657887194403]"
$ws.Range("N14").Value = 0.04357524625438658
$ws.Range("O14").Value = 0.04357524625438658
$ws.Range("Q14").Value = "[-0.8679475199106932, 2.1384214258669267]"
$ws.Range("R14").Value = 0.3991936554758577
$ws.Range("S14").Value = 0.3991936554758577
$ws.Range("U14").Value = "[5.0212774167757255, 12.75369444944139]"
$ws.Range("V14").Value = [double]"3.1191459912705e-05"
$ws.Range("W14").Value = [double]"3.1191459912705e-05"
$ws.Range("Y14").Value = 16.01653653653677
$ws.Range("Z14").Value = 27.6339939939944

# Row 15
$ws.Range("M15").Value = "[0.08636545017252928, 15.207623185822566]"
$ws.Range("N15").Value = 0.047546395746487
$ws.Range("O15").Value = 0.047546395746487
$ws.Range("Q15").Value = "[-1.0943686120613094, 1.7610529389492342]"
$ws.Range("R15").Value = 0.6404453457638981
$ws.Range("S15").Value = 0.6404453457638981
$ws.Range("U15").Value = "[5.783517337109636, 14.140314072566255]"
$ws.Range("V15").Value = [double]"1.773648006309081e-05"
$ws.Range("W15").Value = [double]"1.773648006309081e-05"
$ws.Range("Y15").Value = 17.47479479479505
$ws.Range("Z15").Value = 28.50894894894937
